# Scheduled market-data refresh: update currentAveragePrice / Leve price & profit
# columns (H:N) across all Leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 32
$ws.Range("H32").Value = 1288.625
$ws.Range("J32").Value = 1313.0769
$ws.Range("L32").Value = 1313.0769
$ws.Range("N32").Value = -1965.0769
# row 51
$ws.Range("H51").Value = 52340052
$ws.Range("J51").Value = 4210.5557
$ws.Range("L51").Value = 4210.5557
$ws.Range("N51").Value = -5178.5557
# row 76
$ws.Range("H76").Value = 4248.25
$ws.Range("I76").Value = 4332.3335
$ws.Range("J76").Value = 3996
$ws.Range("K76").Value = 4332.3335
$ws.Range("L76").Value = 3996
$ws.Range("M76").Value = -4017.3335
$ws.Range("N76").Value = -4626
# row 79
$ws.Range("H79").Value = 4248.25
$ws.Range("I79").Value = 4332.3335
$ws.Range("J79").Value = 3996
$ws.Range("K79").Value = 4332.3335
$ws.Range("L79").Value = 3996
$ws.Range("M79").Value = -3240.3335
$ws.Range("N79").Value = -6180
# row 92
$ws.Range("H92").Value = 985.5
$ws.Range("I92").Value = 756.3333
$ws.Range("K92").Value = 756.3333
$ws.Range("M92").Value = 491.6667
# row 106
$ws.Range("H106").Value = 2777.7646
$ws.Range("I106").Value = 1978.7693
$ws.Range("K106").Value = 1978.7693
$ws.Range("M106").Value = -1347.7693
# row 113
$ws.Range("H113").Value = 4553.3335
$ws.Range("I113").Value = 3350
$ws.Range("J113").Value = 4662.727
$ws.Range("K113").Value = 3350
$ws.Range("L113").Value = 4662.727
$ws.Range("M113").Value = -96
$ws.Range("N113").Value = -11170.727
# row 115
$ws.Range("H115").Value = 919.4
$ws.Range("I115").Value = 919.4
$ws.Range("K115").Value = 2758.2
$ws.Range("M115").Value = -1191.2
# row 116
$ws.Range("H116").Value = 3499
$ws.Range("I116").Value = 3288.2222
$ws.Range("K116").Value = 3288.2222
$ws.Range("M116").Value = 153.7777999999998
# row 132
$ws.Range("H132").Value = 40004576
$ws.Range("I132").Value = 47623976
$ws.Range("K132").Value = 142871928
$ws.Range("M132").Value = -142869398
# row 135
$ws.Range("H135").Value = 1698.65
$ws.Range("I135").Value = 1665.2222
$ws.Range("K135").Value = 14986.9998
$ws.Range("M135").Value = -12451.9998
# row 137
$ws.Range("H137").Value = 2501.8484
$ws.Range("I137").Value = 2381.087
$ws.Range("K137").Value = 7143.261
$ws.Range("M137").Value = -4593.261
# row 141
$ws.Range("H141").Value = 5698.4287
$ws.Range("J141").Value = 6397.8
$ws.Range("L141").Value = 19193.4
$ws.Range("N141").Value = -29553.4

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 1144290.2
$ws.Range("I61").Value = 2004816.9
$ws.Range("K61").Value = 2004816.9
$ws.Range("M61").Value = -2004604.9
# row 74
$ws.Range("H74").Value = 4468426
# row 77
$ws.Range("H77").Value = 4468426
# row 110
$ws.Range("H110").Value = 1246.5
$ws.Range("J110").Value = 556.5
$ws.Range("L110").Value = 556.5
$ws.Range("N110").Value = -4646.5
# row 132
$ws.Range("H132").Value = 2080689.8
$ws.Range("I132").Value = 2405394.5
$ws.Range("K132").Value = 7216183.5
$ws.Range("M132").Value = -7213653.5
# row 136
$ws.Range("H136").Value = 1144290.2
$ws.Range("I136").Value = 2004816.9
$ws.Range("K136").Value = 6014450.699999999
$ws.Range("M136").Value = -6011900.699999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 2698.9048
$ws.Range("I20").Value = 2703.7778
$ws.Range("J20").Value = 2669.6667
$ws.Range("K20").Value = 2703.7778
$ws.Range("L20").Value = 2669.6667
$ws.Range("M20").Value = -2456.7778
$ws.Range("N20").Value = -3163.6667
# row 86
$ws.Range("H86").Value = 3514620.2
$ws.Range("I86").Value = 5559006.5
$ws.Range("K86").Value = 5559006.5
$ws.Range("M86").Value = -5557883.5
# row 89
$ws.Range("H89").Value = 3514620.2
$ws.Range("I89").Value = 5559006.5
$ws.Range("K89").Value = 27795032.5
$ws.Range("M89").Value = -27789416.5
# row 105
$ws.Range("H105").Value = 4432.5
$ws.Range("I105").Value = 2993
$ws.Range("J105").Value = 7599.4
$ws.Range("K105").Value = 2993
$ws.Range("L105").Value = 7599.4
$ws.Range("M105").Value = -1246
$ws.Range("N105").Value = -11093.4
# row 134
$ws.Range("H134").Value = 5078.0845
$ws.Range("I134").Value = 1619.8085
$ws.Range("J134").Value = 9593.056
$ws.Range("K134").Value = 4859.4255
$ws.Range("L134").Value = 28779.168
$ws.Range("M134").Value = -2324.4255
$ws.Range("N134").Value = -33849.16800000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2918.2144
$ws.Range("I31").Value = 2624.0908
$ws.Range("J31").Value = 3996.6667
$ws.Range("K31").Value = 2624.0908
$ws.Range("L31").Value = 3996.6667
$ws.Range("M31").Value = -2329.0908
$ws.Range("N31").Value = -4586.6667
# row 34
$ws.Range("H34").Value = 2918.2144
$ws.Range("I34").Value = 2624.0908
$ws.Range("J34").Value = 3996.6667
$ws.Range("K34").Value = 2624.0908
$ws.Range("L34").Value = 3996.6667
$ws.Range("M34").Value = -2422.0908
$ws.Range("N34").Value = -4400.6667
# row 94
$ws.Range("H94").Value = 1241.4
$ws.Range("I94").Value = 1842
$ws.Range("J94").Value = 1091.25
$ws.Range("K94").Value = 1842
$ws.Range("L94").Value = 1091.25
$ws.Range("M94").Value = -1391
$ws.Range("N94").Value = -1993.25
# row 99
$ws.Range("H99").Value = 3176.8572
$ws.Range("I99").Value = 2662
$ws.Range("J99").Value = 4721.4287
$ws.Range("K99").Value = 2662
$ws.Range("L99").Value = 4721.4287
$ws.Range("M99").Value = -1164
$ws.Range("N99").Value = -7717.4287
# row 122
$ws.Range("H122").Value = 2846.6765
$ws.Range("J122").Value = 2532.3333
$ws.Range("L122").Value = 7596.999899999999
$ws.Range("N122").Value = -12496.9999
# row 126
$ws.Range("H126").Value = 3176.8572
$ws.Range("I126").Value = 2662
$ws.Range("J126").Value = 4721.4287
$ws.Range("K126").Value = 7986
$ws.Range("L126").Value = 14164.2861
$ws.Range("M126").Value = -5516
$ws.Range("N126").Value = -19104.2861

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 29
$ws.Range("H29").Value = 3456.8
$ws.Range("J29").Value = 4196.25
$ws.Range("L29").Value = 12588.75
$ws.Range("N29").Value = -13142.75
# row 98
$ws.Range("H98").Value = 2019.8
$ws.Range("I98").Value = 649.5
$ws.Range("K98").Value = 1948.5
$ws.Range("M98").Value = -450.5
# row 107
$ws.Range("H107").Value = 2961.9167
$ws.Range("I107").Value = 261.5
$ws.Range("K107").Value = 784.5
$ws.Range("M107").Value = 1135.5
# row 122
$ws.Range("H122").Value = 1135.6666
$ws.Range("I122").Value = 682
$ws.Range("J122").Value = 1438.1111
$ws.Range("K122").Value = 6138
$ws.Range("L122").Value = 12942.9999
$ws.Range("M122").Value = -3688
$ws.Range("N122").Value = -17842.9999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 15360.641
$ws.Range("I102").Value = 2293.742
$ws.Range("J102").Value = 65994.875
$ws.Range("K102").Value = 2293.742
$ws.Range("L102").Value = 65994.875
$ws.Range("M102").Value = -671.7420000000002
$ws.Range("N102").Value = -69238.875
# row 121
$ws.Range("H121").Value = 89458
$ws.Range("J121").Value = 89458
$ws.Range("L121").Value = 89458
$ws.Range("N121").Value = -92952
# row 122
$ws.Range("H122").Value = 6392.7144
$ws.Range("I122").Value = 5428.857
$ws.Range("J122").Value = 7356.5713
$ws.Range("K122").Value = 16286.571
$ws.Range("L122").Value = 22069.7139
$ws.Range("M122").Value = -13836.571
$ws.Range("N122").Value = -26969.7139

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value = 295.77777
$ws.Range("I55").Value = 294.83334
$ws.Range("K55").Value = 294.83334
$ws.Range("M55").Value = -121.83334
# row 93
$ws.Range("H93").Value = 1450.5333
$ws.Range("I93").Value = 1309.8889
$ws.Range("J93").Value = 1661.5
$ws.Range("K93").Value = 1309.8889
$ws.Range("L93").Value = 1661.5
$ws.Range("M93").Value = -61.88889999999992
$ws.Range("N93").Value = -4157.5
# row 137
$ws.Range("H137").Value = 30000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# row 139
$ws.Range("H139").Value = 57538.08
$ws.Range("J139").Value = 62545
$ws.Range("L139").Value = 62545
$ws.Range("N139").Value = -72825
# row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = 0

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 96
$ws.Range("H96").Value = 3942.077
$ws.Range("J96").Value = 4114.9
$ws.Range("L96").Value = 4114.9
$ws.Range("N96").Value = -6860.9
# row 100
$ws.Range("H100").Value = 881.6070999999999
$ws.Range("I100").Value = 780.6799999999999
$ws.Range("K100").Value = 1561.36
$ws.Range("M100").Value = -1020.36
# row 122
$ws.Range("H122").Value = 2583
$ws.Range("I122").Value = 2274.1904
$ws.Range("J122").Value = 3663.8333
$ws.Range("K122").Value = 6822.5712
$ws.Range("L122").Value = 10991.4999
$ws.Range("M122").Value = -4372.5712
$ws.Range("N122").Value = -15891.4999
# row 126
$ws.Range("H126").Value = 1656.16
$ws.Range("I126").Value = 1197.579
$ws.Range("K126").Value = 3592.737
$ws.Range("M126").Value = -1122.737
# row 136
$ws.Range("H136").Value = 6466.4634
$ws.Range("I136").Value = 6746.8203
$ws.Range("K136").Value = 20240.4609
$ws.Range("M136").Value = -17690.4609
